$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.39907431602478
$ws.Range("B1").Value = 1.430711030960083
$ws.Range("C1").Value = 1.537987589836121
$ws.Range("D1").Value = 2.202733278274536
$ws.Range("E1").Value = 4.348400592803955
